$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text for the two existing "TASK(HTML/CSS)" rows:
# the mock-page task text is replaced with the template task text.
$ws.Range("C9").Value = "To Create template of Amazon.in(only homepage)"
$ws.Range("C10").Value = "To Create template of Amazon.in(only homepage)"

# Add a new tracker row (row 12) for a new JavaScript topic, copying the
# date-cell number formatting from the row above so the style matches.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 43321

$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = "JavaScript"

$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = "Variables,operators,String Interpolation,Control Flow,Functions,"

[void]$ws.Range("C12").Select()
